$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: flip S311:S318 from 0 to 1 ---
$ws.Cells.Item(311, 19).Value = 1
$ws.Cells.Item(312, 19).Value = 1
$ws.Cells.Item(313, 19).Value = 1
$ws.Cells.Item(314, 19).Value = 1
$ws.Cells.Item(315, 19).Value = 1
$ws.Cells.Item(316, 19).Value = 1
$ws.Cells.Item(317, 19).Value = 1
$ws.Cells.Item(318, 19).Value = 1

# --- Step 2: add new cells / rows, in an order that reproduces the exact
#     shared-string first-use sequence from the target workbook ---
$ws.Cells.Item(319, 3).Value = "maneti2018"
$ws.Cells.Item(318, 26).Value = "exposure environment"
$ws.Cells.Item(319, 26).Value = "exposure environment"
$ws.Cells.Item(318, 27).Value = "fluctuating"
$ws.Cells.Item(319, 27).Value = "constant"
$ws.Cells.Item(319, 11).Value = "Time to heat knockdown"
$ws.Cells.Item(319, 12).Value = "min"
$ws.Cells.Item(319, 22).Value = "simulans"
$ws.Cells.Item(324, 3).Value = "peng2014"
$ws.Cells.Item(57, 4).Copy()
$ws.Cells.Item(324, 4).PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(324, 4).Value = "Figure 1a"
$ws.Cells.Item(57, 4).Copy()
$ws.Cells.Item(325, 4).PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(325, 4).Value = "Figure 1b"
$ws.Cells.Item(57, 4).Copy()
$ws.Cells.Item(326, 4).PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(326, 4).Value = "Figure 1c"
$ws.Cells.Item(57, 4).Copy()
$ws.Cells.Item(327, 4).PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(327, 4).Value = "Figure 1d"
$ws.Cells.Item(326, 11).Value = "LTmax"
$ws.Cells.Item(327, 11).Value = "LTmin"
$ws.Cells.Item(324, 21).Value = "Parabramis "
$ws.Cells.Item(324, 22).Value = "pekinensis"
$ws.Cells.Item(319, 1).Value = 1
$ws.Cells.Item(319, 2).Value = 1
$ws.Cells.Item(319, 4).Value = "Figure 1"
$ws.Cells.Item(319, 5).Value = 0
$ws.Cells.Item(319, 6).Value = 1
$ws.Cells.Item(318, 7).Copy()
$ws.Cells.Item(319, 7).PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(319, 7).Value = 37.5
$ws.Cells.Item(319, 8).Value = 23
$ws.Cells.Item(319, 9).Value = 15
$ws.Cells.Item(319, 10).Value = 24
$ws.Cells.Item(151, 13).Copy()
$ws.Cells.Item(319, 13).PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(319, 13).Value = 61.776134300000002
$ws.Cells.Item(151, 13).Copy()
$ws.Cells.Item(319, 14).PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(319, 14).Value = 58.528296699999999
$ws.Cells.Item(319, 15).Value = 1.7616221999999999
$ws.Cells.Item(319, 16).Value = 2.2107081200000001
$ws.Cells.Item(319, 17).Value = 192
$ws.Cells.Item(319, 18).Value = 192
$ws.Cells.Item(319, 19).Value = 1
$ws.Cells.Item(319, 20).Value = 0
$ws.Cells.Item(319, 21).Value = "Drosophila "
$ws.Cells.Item(319, 23).Value = 1
$ws.Cells.Item(319, 24).Value = 2
$ws.Cells.Item(319, 25).Value = 1
$ws.Cells.Item(320, 1).Value = 1
$ws.Cells.Item(320, 2).Value = 2
$ws.Cells.Item(320, 3).Value = "maneti2018"
$ws.Cells.Item(320, 4).Value = "Figure 1"
$ws.Cells.Item(320, 5).Value = 0
$ws.Cells.Item(320, 6).Value = 1
$ws.Cells.Item(318, 7).Copy()
$ws.Cells.Item(320, 7).PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(320, 7).Value = 37.5
$ws.Cells.Item(320, 8).Value = 23
$ws.Cells.Item(320, 9).Value = 15
$ws.Cells.Item(320, 10).Value = 24
$ws.Cells.Item(320, 11).Value = "Time to heat knockdown"
$ws.Cells.Item(320, 12).Value = "min"
$ws.Cells.Item(151, 13).Copy()
$ws.Cells.Item(320, 13).PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(320, 13).Value = 36.0766524
$ws.Cells.Item(151, 13).Copy()
$ws.Cells.Item(320, 14).PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(320, 14).Value = 29.7198633
$ws.Cells.Item(320, 15).Value = 1.8307426600000001
$ws.Cells.Item(320, 16).Value = 1.51986183
$ws.Cells.Item(320, 17).Value = 192
$ws.Cells.Item(320, 18).Value = 192
$ws.Cells.Item(320, 19).Value = 1
$ws.Cells.Item(320, 20).Value = 0
$ws.Cells.Item(320, 21).Value = "Drosophila "
$ws.Cells.Item(320, 22).Value = "simulans"
$ws.Cells.Item(320, 23).Value = 1
$ws.Cells.Item(320, 24).Value = 2
$ws.Cells.Item(320, 25).Value = 1
$ws.Cells.Item(321, 1).Value = 1
$ws.Cells.Item(321, 2).Value = 3
$ws.Cells.Item(321, 3).Value = "maneti2018"
$ws.Cells.Item(321, 4).Value = "Figure 2"
$ws.Cells.Item(321, 5).Value = 0
$ws.Cells.Item(321, 6).Value = 1
$ws.Cells.Item(321, 8).Value = 23
$ws.Cells.Item(321, 9).Value = 15
$ws.Cells.Item(321, 10).Value = 24
$ws.Cells.Item(321, 11).Value = "CTmax"
$ws.Cells.Item(321, 12).Value = "min"
$ws.Cells.Item(321, 13).Value = 39.054339116499399
$ws.Cells.Item(321, 14).Value = 40.230621169999999
$ws.Cells.Item(151, 13).Copy()
$ws.Cells.Item(321, 15).PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(321, 15).Value = 0.11858974
$ws.Cells.Item(321, 16).Value = 0.096153846153850253
$ws.Cells.Item(321, 17).Value = 16
$ws.Cells.Item(321, 18).Value = 16
$ws.Cells.Item(321, 19).Value = 1
$ws.Cells.Item(321, 20).Value = 0
$ws.Cells.Item(321, 21).Value = "Drosophila "
$ws.Cells.Item(321, 22).Value = "simulans"
$ws.Cells.Item(321, 23).Value = 1
$ws.Cells.Item(321, 24).Value = 2
$ws.Cells.Item(321, 25).Value = 1
$ws.Cells.Item(322, 1).Value = 1
$ws.Cells.Item(322, 2).Value = 3
$ws.Cells.Item(322, 3).Value = "maneti2018"
$ws.Cells.Item(322, 4).Value = "Figure 2"
$ws.Cells.Item(322, 5).Value = 0
$ws.Cells.Item(322, 6).Value = 1
$ws.Cells.Item(322, 8).Value = 23
$ws.Cells.Item(322, 9).Value = 15
$ws.Cells.Item(322, 10).Value = 24
$ws.Cells.Item(322, 11).Value = "CTmax"
$ws.Cells.Item(322, 12).Value = "min"
$ws.Cells.Item(322, 13).Value = 39.445207984949803
$ws.Cells.Item(322, 14).Value = 39.852259267001102
$ws.Cells.Item(322, 15).Value = 0.10576923076919797
$ws.Cells.Item(322, 16).Value = 0.081730769230752287
$ws.Cells.Item(322, 17).Value = 16
$ws.Cells.Item(322, 18).Value = 16
$ws.Cells.Item(322, 19).Value = 1
$ws.Cells.Item(322, 20).Value = 0
$ws.Cells.Item(322, 21).Value = "Drosophila "
$ws.Cells.Item(322, 22).Value = "simulans"
$ws.Cells.Item(322, 23).Value = 1
$ws.Cells.Item(322, 24).Value = 2
$ws.Cells.Item(322, 25).Value = 1
$ws.Cells.Item(323, 1).Value = 1
$ws.Cells.Item(323, 2).Value = 3
$ws.Cells.Item(323, 3).Value = "maneti2018"
$ws.Cells.Item(323, 4).Value = "Figure 2"
$ws.Cells.Item(323, 5).Value = 0
$ws.Cells.Item(323, 6).Value = 1
$ws.Cells.Item(323, 8).Value = 23
$ws.Cells.Item(323, 9).Value = 15
$ws.Cells.Item(323, 10).Value = 24
$ws.Cells.Item(323, 11).Value = "CTmax"
$ws.Cells.Item(323, 12).Value = "min"
$ws.Cells.Item(323, 13).Value = 39.637358904682202
$ws.Cells.Item(323, 14).Value = 40.005948648272003
$ws.Cells.Item(323, 15).Value = 0.11858974358970187
$ws.Cells.Item(323, 16).Value = 0.068914611204000664
$ws.Cells.Item(323, 17).Value = 16
$ws.Cells.Item(323, 18).Value = 16
$ws.Cells.Item(323, 19).Value = 1
$ws.Cells.Item(323, 20).Value = 0
$ws.Cells.Item(323, 21).Value = "Drosophila "
$ws.Cells.Item(323, 22).Value = "simulans"
$ws.Cells.Item(323, 23).Value = 1
$ws.Cells.Item(323, 24).Value = 2
$ws.Cells.Item(323, 25).Value = 1
$ws.Cells.Item(324, 1).Value = 1
$ws.Cells.Item(324, 2).Value = 1
$ws.Cells.Item(324, 5).Value = 0
$ws.Cells.Item(324, 6).Value = 1
$ws.Cells.Item(324, 8).Value = 20
$ws.Cells.Item(324, 9).Value = 10
$ws.Cells.Item(324, 10).Value = 24
$ws.Cells.Item(324, 11).Value = "CTmax"
$ws.Cells.Item(324, 12).Value = "C"
$ws.Cells.Item(151, 13).Copy()
$ws.Cells.Item(324, 13).PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(324, 13).Value = 33.283018900000002
$ws.Cells.Item(151, 13).Copy()
$ws.Cells.Item(324, 14).PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(324, 14).Value = 35.264150899999997
$ws.Cells.Item(324, 15).Value = 0.79245283
$ws.Cells.Item(324, 16).Value = 0.59433961999999996
$ws.Cells.Item(324, 17).Value = 8
$ws.Cells.Item(324, 18).Value = 8
$ws.Cells.Item(324, 19).Value = 1
$ws.Cells.Item(324, 20).Value = 0
$ws.Cells.Item(324, 23).Value = 1
$ws.Cells.Item(324, 24).Value = 1
$ws.Cells.Item(324, 25).Value = 2
$ws.Cells.Item(325, 1).Value = 2
$ws.Cells.Item(325, 2).Value = 1
$ws.Cells.Item(325, 3).Value = "peng2014"
$ws.Cells.Item(325, 5).Value = 0
$ws.Cells.Item(325, 6).Value = 1
$ws.Cells.Item(325, 8).Value = 20
$ws.Cells.Item(325, 9).Value = 10
$ws.Cells.Item(325, 10).Value = 24
$ws.Cells.Item(325, 11).Value = "CTmin"
$ws.Cells.Item(325, 12).Value = "C"
$ws.Cells.Item(151, 13).Copy()
$ws.Cells.Item(325, 13).PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(325, 13).Value = 5.1844155799999996
$ws.Cells.Item(151, 13).Copy()
$ws.Cells.Item(325, 14).PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(325, 14).Value = 4.8772517799999999
$ws.Cells.Item(325, 15).Value = 0.25345622000000001
$ws.Cells.Item(325, 16).Value = 0.18433179999999999
$ws.Cells.Item(325, 17).Value = 8
$ws.Cells.Item(325, 18).Value = 8
$ws.Cells.Item(325, 19).Value = 1
$ws.Cells.Item(325, 20).Value = 0
$ws.Cells.Item(325, 21).Value = "Parabramis "
$ws.Cells.Item(325, 22).Value = "pekinensis"
$ws.Cells.Item(325, 23).Value = 1
$ws.Cells.Item(325, 24).Value = 1
$ws.Cells.Item(325, 25).Value = 2
$ws.Cells.Item(326, 1).Value = 3
$ws.Cells.Item(326, 2).Value = 1
$ws.Cells.Item(326, 3).Value = "peng2014"
$ws.Cells.Item(326, 5).Value = 0
$ws.Cells.Item(326, 6).Value = 1
$ws.Cells.Item(326, 8).Value = 20
$ws.Cells.Item(326, 9).Value = 10
$ws.Cells.Item(326, 10).Value = 24
$ws.Cells.Item(326, 12).Value = "C"
$ws.Cells.Item(151, 13).Copy()
$ws.Cells.Item(326, 13).PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(326, 13).Value = 34.194690299999998
$ws.Cells.Item(326, 14).Value = 36.053097299999997
$ws.Cells.Item(326, 15).Value = 0.37168141999999998
$ws.Cells.Item(326, 16).Value = 0.18584070999999999
$ws.Cells.Item(326, 17).Value = 8
$ws.Cells.Item(326, 18).Value = 8
$ws.Cells.Item(326, 19).Value = 1
$ws.Cells.Item(326, 20).Value = 0
$ws.Cells.Item(326, 21).Value = "Parabramis "
$ws.Cells.Item(326, 22).Value = "pekinensis"
$ws.Cells.Item(326, 23).Value = 1
$ws.Cells.Item(326, 24).Value = 1
$ws.Cells.Item(326, 25).Value = 2
$ws.Cells.Item(327, 1).Value = 4
$ws.Cells.Item(327, 2).Value = 1
$ws.Cells.Item(327, 3).Value = "peng2014"
$ws.Cells.Item(327, 5).Value = 0
$ws.Cells.Item(327, 6).Value = 1
$ws.Cells.Item(327, 8).Value = 20
$ws.Cells.Item(327, 9).Value = 10
$ws.Cells.Item(327, 10).Value = 24
$ws.Cells.Item(327, 12).Value = "C"
$ws.Cells.Item(151, 13).Copy()
$ws.Cells.Item(327, 13).PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(327, 13).Value = 3.4199134199999999
$ws.Cells.Item(151, 13).Copy()
$ws.Cells.Item(327, 14).PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(327, 14).Value = 2.5541125500000001
$ws.Cells.Item(327, 15).Value = 0.19480518999999999
$ws.Cells.Item(327, 16).Value = 0.17316017
$ws.Cells.Item(327, 17).Value = 8
$ws.Cells.Item(327, 18).Value = 8
$ws.Cells.Item(327, 19).Value = 1
$ws.Cells.Item(327, 20).Value = 0
$ws.Cells.Item(327, 21).Value = "Parabramis "
$ws.Cells.Item(327, 22).Value = "pekinensis"
$ws.Cells.Item(327, 23).Value = 1
$ws.Cells.Item(327, 24).Value = 1
$ws.Cells.Item(327, 25).Value = 2
